$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999998593648143,
    0.99531893239156999,
    0.98144903717425758,
    0.97521877446479577,
    0.96601031874524734,
    0.94421381987814001,
    0.94469903265270216,
    0.94362476433294651,
    0.94702190820974164,
    0.95140410476153014,
    0.95201357679876042,
    0.95212125043906104,
    0.95423107512129035,
    0.95588857584347187,
    0.95885858986718031,
    0.96300448854724097,
    0.97326398684909488,
    0.97883816579941518,
    0.99312447520142677,
    0.98600787155506198,
    0.98460946664570703,
    0.98334498472844167,
    0.97726710178208909,
    0.97061927872649512,
    0.96416240832850564,
    0.94055899724449232,
    0.93571312125714057,
    0.91424254279953709,
    0.89897242190751991,
    0.89240244695079318,
    0.88474937441165413,
    0.88307018771941825,
    0.88255021291278068
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
